# Map behavior/call_flow/version impact summary
#
# Inserts three "[onshow.*]" placeholder paragraphs under the
# "Functional Behavior" and "Call flow" bullets, and adds a brand new
# "Version Impact Summary" bullet (with its own "[onshow.summary]"
# placeholder) right after "Call flow".
#
# Strategy: Range.InsertXML() *replaces* the contents of the range it
# is called on, so existing paragraphs are rewritten in-place (keeping
# their own <w:p> element, only with the runs/pPr we specify) by
# feeding InsertXML the full desired <w:p>...</w:p> fragment for that
# paragraph. New paragraphs are created first via
# Range.InsertParagraphAfter() (which yields an empty <w:p/> with no
# legacy rsid attributes right after the anchor paragraph) and then
# their content is filled in the same way.
#
# Paragraph.Index is unreliable in this document (it has tables, and
# Index appears to be scoped per container), so paragraphs are always
# located/re-located by their 1-based position in $d.Paragraphs(...).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$ns = $wNs + ' ' + $w14Ns

function Get-ParagraphPosByText($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13)
        $t = $t.TrimEnd()
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# --- "Functional Behavior" bullet: append two trailing space runs ---
$posBehaviorHeading = Get-ParagraphPosByText("Functional Behavior")
$xmlBehaviorHeading = '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="15"/></w:numPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Functional Behavior</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posBehaviorHeading).Range.InsertXML($xmlBehaviorHeading)

# --- New paragraph: "[onshow.behavior]" right after "Functional Behavior" ---
$d.Paragraphs($posBehaviorHeading).Range.InsertParagraphAfter()
$posBehaviorTag = $posBehaviorHeading + 1
$xmlBehaviorTag = '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Normal"/><w:ind w:left="2124" w:firstLine="708"/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>onshow</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.behavior</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>]</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posBehaviorTag).Range.InsertXML($xmlBehaviorTag)

# --- "Call flow" bullet: mark w14:noSpellErr="1" ---
$posCallFlowHeading = Get-ParagraphPosByText("Call flow")
$xmlCallFlowHeading = '<w:p ' + $ns + ' w14:noSpellErr="1">' +
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="15"/></w:numPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Call flow</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posCallFlowHeading).Range.InsertXML($xmlCallFlowHeading)

# --- New paragraph: "[onshow.call_flow]" right after "Call flow" ---
$d.Paragraphs($posCallFlowHeading).Range.InsertParagraphAfter()
$posCallFlowTag = $posCallFlowHeading + 1
$xmlCallFlowTag = '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Normal"/><w:ind w:left="2124" w:firstLine="708"/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[onshow</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.call_flow</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>]</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posCallFlowTag).Range.InsertXML($xmlCallFlowTag)

# --- New bullet "Version Impact Summary" right after "[onshow.call_flow]" ---
$d.Paragraphs($posCallFlowTag).Range.InsertParagraphAfter()
$posVersionHeading = $posCallFlowTag + 1
$xmlVersionHeading = '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="15"/></w:numPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Version Impact Summary</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posVersionHeading).Range.InsertXML($xmlVersionHeading)

# --- New paragraph: "[onshow.summary]" right after "Version Impact Summary" ---
$d.Paragraphs($posVersionHeading).Range.InsertParagraphAfter()
$posVersionTag = $posVersionHeading + 1
$xmlVersionTag = '<w:p ' + $ns + '>' +
    '<w:pPr><w:pStyle w:val="Normal"/><w:ind w:left="2124" w:firstLine="708"/></w:pPr>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>[onshow</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.summary</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>]</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs($posVersionTag).Range.InsertXML($xmlVersionTag)

Write-Output "done"
